# Apply the "new routine" update to the L5CG9 3rd semester routine sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: collapse A1:L1 into a single title cell in A1 ---
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value2 = "Herald College Kathmandu"

# --- Data rows 2-10: re-derive the new column layout from the old one ---
# Old layout: A Day | B Time | C Module Code | D Module Title | E Hours |
#             F Class Type | G Lecturer | H Room | I Block | J Group | K Level | L Course
# New layout: A Day | B Time | C Hours | D Module Code | E Module Title |
#             F Class Type | G Lecturer | H Group | I Block | J Room
for ($r = 2; $r -le 10; $r++) {
    # Read every old value first (via Value2, which returns the real scalar
    # instead of a property-descriptor string) so later writes never clobber
    # a still-needed source.
    $moduleCode  = $ws.Cells.Item($r, 3).Value2   # C (old Module Code)
    $moduleTitle = $ws.Cells.Item($r, 4).Value2   # D (old Module Title)
    $hours       = $ws.Cells.Item($r, 5).Value2   # E (old Hours)
    $room        = $ws.Cells.Item($r, 8).Value2   # H (old Room)
    $group       = $ws.Cells.Item($r, 10).Value2  # J (old Group)

    # Now write the new values into their new positions.
    $ws.Cells.Item($r, 3).Value2 = $hours         # C: Hours (numeric)
    $ws.Cells.Item($r, 4).Value2 = $moduleCode    # D: Module Code
    $ws.Cells.Item($r, 5).Value2 = $moduleTitle   # E: Module Title
    $ws.Cells.Item($r, 8).Value2 = $group         # H: Group
    $ws.Cells.Item($r, 10).Value2 = $room         # J: Room (I already holds Block = WLV/HCK, unchanged)

    # Clear the now-unused trailing columns (old Level, Course)
    $ws.Cells.Item($r, 11).ClearContents()
    $ws.Cells.Item($r, 12).ClearContents()
}

# --- Shrink the used dimension to match the new layout (A1:J10) ---
$ws.Range("K1:L10").ClearContents()
